$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 55 (pushes existing rows 55.. down by one)
$ws.Rows.Item(55).Insert()
$ws.Range("A55").Value = 11
$ws.Range("B55").Value = "Vega Monumental Concepción"
$ws.Range("C55").Value = "Bíobío"
$ws.Range("D55").Value = 44645
$ws.Range("E55").Value = 8
$ws.Range("F55").Value = 100112003
$ws.Range("G55").Value = "Ajo"
$ws.Range("H55").Value = "Chino"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 220
$ws.Range("K55").Value = 18000
$ws.Range("L55").Value = 19000
$ws.Range("M55").Value = 18545
$ws.Range("N55").Value = "$/caja 10 kilos"
$ws.Range("O55").Value = "China"
$ws.Range("P55").Value = 1854
$ws.Range("Q55").Value = 10
$ws.Range("R55").Value = "Hortaliza"

# Insert a second new data row at row 131 (pushes existing rows 131.. down by one)
$ws.Rows.Item(131).Insert()
$ws.Range("A131").Value = 11
$ws.Range("B131").Value = "Vega Monumental Concepción"
$ws.Range("C131").Value = "Bíobío"
$ws.Range("D131").Value = 44644
$ws.Range("E131").Value = 8
$ws.Range("F131").Value = 100112003
$ws.Range("G131").Value = "Ajo"
$ws.Range("H131").Value = "Chino"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 170
$ws.Range("K131").Value = 19000
$ws.Range("L131").Value = 20000
$ws.Range("M131").Value = 19529
$ws.Range("N131").Value = "$/caja 10 kilos"
$ws.Range("O131").Value = "China"
$ws.Range("P131").Value = 1953
$ws.Range("Q131").Value = 10
$ws.Range("R131").Value = "Hortaliza"
